$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.892.11'
$ws.Range("E2").Value = '  -1.62%  '
$ws.Range("D3").Value = '3.304.15'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = '557.28'
$ws.Range("E5").Value = '  -0.52%  '
$fmt = $ws.Range("D6").NumberFormat
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '185.40'
$ws.Range("D6").NumberFormat = $fmt
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E8").Value = '  -2.19%  '
$ws.Range("D9").Value = '3.296.19'
$ws.Range("E9").Value = '  -1.99%  '
$ws.Range("E10").Value = '  -3.48%  '
$ws.Range("E11").Value = '  -2.69%  '
$ws.Range("D12").Value = '47.45'
$ws.Range("E12").Value = '  -1.40%  '
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("E14").Value = '  -1.82%  '
$ws.Range("D15").Value = '629.35'
$ws.Range("E15").Value = '  +2.73%  '
$ws.Range("D16").Value = '3.834.64'
$ws.Range("E16").Value = '  -1.86%  '
$fmt = $ws.Range("D17").NumberFormat
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.10'
$ws.Range("D17").NumberFormat = $fmt
$ws.Range("E17").Value = '  +2.17%  '
$ws.Range("D18").Value = '65.905.29'
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("E19").Value = '  -1.76%  '
$ws.Range("D20").Value = '3.295.80'
$ws.Range("E20").Value = '  -3.53%  '
$ws.Range("D21").Value = '11.34'
$ws.Range("E21").Value = '  -3.55%  '
$ws.Range("E22").Value = '  -1.35%  '
$ws.Range("D23").Value = '18.03'
$ws.Range("E23").Value = '  +3.77%  '
$ws.Range("D24").Value = '102.48'
$ws.Range("E24").Value = '  +7.64%  '
$ws.Range("D25").Value = '4.97'
$ws.Range("E25").Value = '  -3.05%  '
$ws.Range("E26").Value = '  -4.22%  '
$ws.Range("D27").Value = '5.95'
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("E28").Value = '  -1.85%  '
$ws.Range("E29").Value = '  -1.06%  '
$ws.Range("D30").Value = '8.64'
$ws.Range("E30").Value = '  -2.50%  '
$ws.Range("D31").Value = '30.19'
$ws.Range("E31").Value = '  -2.41%  '
$ws.Range("D32").Value = '4.05'
$ws.Range("E32").Value = '  +3.34%  '
$ws.Range("D33").Value = '6.39'
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("E34").Value = '  -1.59%  '
$ws.Range("D35").Value = '547.97'
$ws.Range("E35").Value = '  -1.26%  '
$ws.Range("E36").Value = '  -0.66%  '
$ws.Range("D37").Value = '3.797.43'
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").Value = '57.44'
$ws.Range("E38").Value = '  -2.79%  '
$ws.Range("E39").Value = '  +0.42%  '
$ws.Range("E40").Value = '  +0.86%  '
$ws.Range("D41").Value = '33.59'
$ws.Range("E41").Value = '  +2.45%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("E48").Value = '  -1.53%  '
$ws.Range("E49").Value = '  -1.90%  '
$ws.Range("E50").Value = '  -4.53%  '
$ws.Range("E51").Value = '  -0.09%  '

# Rows 43-47: coin list reordered (ranking shuffled) with updated price/volume
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '3.24'
$ws.Range("E43").Value = '  -7.08%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").Value = '2.69'
$ws.Range("E44").Value = '  -2.24%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").Value = '0.334'
$ws.Range("E45").Value = '  -6.08%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '3.27'
$ws.Range("E46").Value = '  +0.62%  '
$ws.Range("B47").Value = 'CoreDAO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D47").Value = '3.16'
$ws.Range("E47").Value = '  -16.08%  '
